$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# B11 needs to hold the literal text "1" (a shared string), not the number 1,
# while keeping its existing cell style (s="23") untouched. A plain
# Range.Value assignment of a numeric-looking string gets auto-coerced to a
# number by Excel (same as typing 1 into a General-formatted cell), and
# forcing text via NumberFormat="@" / a leading apostrophe creates a brand
# new cell style. Instead, stage the text in a scratch cell (via a formula
# that evaluates to the text "1"), copy it, and paste-special just the value
# into B11 - this preserves B11's existing style while writing a genuine
# text cell.
$scratch = $ws.Range("B1")
$scratch.Formula = "=""1"""
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

